# Update AuditRowUpdaterServiceTests for new UT config and profile.
#
# 1. Rename several vendor sheets to their shortened display names.
# 2. Move the "active" sheet/tab from Fortinet Technologies to SAP
#    (clears tabSelected on the old tab, sets it on the new one, and
#    updates the workbook's activeTab index).
# 3. Update the remembered cell selection on a few sheets.

$wb = $excel.ActiveWorkbook

# --- Rename sheets -------------------------------------------------------
$wb.Worksheets.Item("Cisco Systems, Inc.").Name = "CISCO"
$wb.Worksheets.Item("Fortinet Technologies Inc").Name = "Fortinet Technologies"
$wb.Worksheets.Item("Oracle Corporation").Name = "Oracle"
$wb.Worksheets.Item("Waves Audio Ltd.").Name = "Waves Audio"

# --- Update remembered selections ----------------------------------------
$wsVendorNotFound = $wb.Worksheets.Item("Vendor Not Found")
$wsVendorNotFound.Activate()
$wsVendorNotFound.Range("C11").Select() | Out-Null

$wsCisco = $wb.Worksheets.Item("CISCO")
$wsCisco.Activate()
$wsCisco.Range("B14").Select() | Out-Null

$wsFortinet = $wb.Worksheets.Item("Fortinet Technologies")
$wsFortinet.Activate()
$wsFortinet.Range("C4").Select() | Out-Null

# --- Move the active tab to SAP and update its selection -----------------
$wsSAP = $wb.Worksheets.Item("SAP")
$wsSAP.Activate()
$wsSAP.Range("E5").Select() | Out-Null
